# Update the "Förändrad" date column (C) for rows 2-33 from 45243 to 45244
# (i.e. change the stored serial date value by +1 day, 2023-11-13 -> 2023-11-14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
